# Apply edits described by the commit "Fixed some small stuff and the itembank"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "Fear" item labels to "Fearful" in column B (rows 7-11)
$ws.Range("B7").Value = "adagio20Fearful2"
$ws.Range("B8").Value = "adagio07Fearful"
$ws.Range("B9").Value = "adagio16Fearful"
$ws.Range("B10").Value = "adagio22Fearful"
$ws.Range("B11").Value = "adagio13Fearful"

# Widen column B so the longer labels are fully visible
$ws.Columns.Item(2).ColumnWidth = 47.67

# Update the active selection on the sheet
$ws.Range("C20").Select()
